$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 186, shifting rows 186:266 down to 187:267
$ws.Rows.Item(186).Insert()

# Populate the newly inserted row 186 with its values
$ws.Range("A186").Value = 7
$ws.Range("B186").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C186").Value = "Ñuble"
$ws.Range("D186").Value = 44917
$ws.Range("E186").Value = 16
$ws.Range("F186").Value = 100112032
$ws.Range("G186").Value = "Zapallo italiano"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 300
$ws.Range("K186").Value = 4500
$ws.Range("L186").Value = 5000
$ws.Range("M186").Value = 4750
$ws.Range("N186").Value = "$/caja 50 unidades"
$ws.Range("O186").Value = "Región del Maule"
$ws.Range("P186").Value = 95
$ws.Range("Q186").Value = 50
$ws.Range("R186").Value = "Hortaliza"
